$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# The original workbook had Sheet1/Sheet2/Sheet3 all empty. The edit keeps
# only one sheet (the populated "test cases" sheet), named "Sheet1", and
# drops the other two. Deleting Sheet1 & Sheet2 and renaming Sheet3 is what
# reproduces the target sheetId (3) and r:id (rId1) for the surviving sheet.
$null = $wb.Worksheets.Item("Sheet1").Delete()
$null = $wb.Worksheets.Item("Sheet2").Delete()
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Name = "Sheet1"

# Fill in the test-case labels. Writing TC_1..TC_10 before "TestCases"
# reproduces the shared-string table order of the target file (TC_1..TC_10,
# then TestCases last).
for ($i = 1; $i -le 10; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = "TC_$i"
}
$ws.Cells.Item(1, 1).Value = "TestCases"

# Row heights: rows 1-9 use a tall 75pt row (for wrapped/large text), row 10
# uses 50.1pt, row 11 keeps the default.
for ($r = 1; $r -le 9; $r++) {
    $ws.Rows.Item($r).RowHeight = 75
}
$ws.Rows.Item(10).RowHeight = 50.1

# Column B is widened (used for a second column of notes next to the list).
$ws.Columns.Item(2).ColumnWidth = 31.66

# Final selection lands just past the last row of data.
$null = $ws.Range("A12").Select()
